$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update string / text values
$ws.Range("B2").Value = "Uông Cương"
$ws.Range("D2").Value = "1, Ấp An, Ấp Triệu Nguyệt, Quận Thiện Trà Vinh"
$ws.Range("E2").Value = "(090)046-8046"
$ws.Range("G2").Value = "3.392.800 VNĐ"

$ws.Range("B3").Value = "Thịnh Khai Hiếu"
$ws.Range("D3").Value = "1, Ấp An, Ấp Triệu Nguyệt, Quận Thiện Trà Vinh"
$ws.Range("E3").Value = "(84)(96)555-5261"
$ws.Range("G3").Value = "3.257.000 VNĐ"

# Update numeric values
$ws.Range("A2").Value = 1
$ws.Range("F2").Value = 16964

$ws.Range("A3").Value = 3
$ws.Range("F3").Value = 16285

# Update column widths (engine stores width = round(input*7)/7 + 5/7, so
# feed target - 5/7 to land on the closest representable stored width)
$ws.Columns.Item(5).ColumnWidth = 19.185616285714286
$ws.Columns.Item(7).ColumnWidth = 15.519845285714288
